$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded "Before"/"After" columns (and the "Change" column's
# numeric data), shifting everything left. This also drops the no-longer
# referenced shared strings and shrinks the used range to A1:F3.
$ws.Range("G1:H3").Delete(-4161) | Out-Null

# Repurpose column F as a required "Group" column.
$ws.Range("F1").Value2 = "Group"
$ws.Range("F2").Value2 = "A"
$ws.Range("F3").ClearContents() | Out-Null

# Mark the data rows as explicitly carrying the (now required) Normal style,
# which materializes a dedicated cell format for A2:E3.
$ws.Range("A2:E3").Style = "Normal"

$ws.Range("F2").Select() | Out-Null
